$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.614.84'
$ws.Range("E2").Value = '  -7.17%  '
$ws.Range("D3").Value = '1.695.32'
$ws.Range("E3").Value = '  -5.63%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.20'
$ws.Range("E5").Value = '  -4.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5122'
$ws.Range("E6").Value = '  -12.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.006'
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2656'
$ws.Range("E8").Value = '  -4.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '22.19'
$ws.Range("E9").Value = '  -4.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06337'
$ws.Range("E10").Value = '  -6.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07366'
$ws.Range("E11").Value = '  -2.21%  '
$ws.Range("D12").Value = '1.695.67'
$ws.Range("E12").Value = '  -5.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.531'
$ws.Range("E13").Value = '  -5.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5796'
$ws.Range("E14").Value = '  -5.87%  '
$ws.Range("D15").Value = '1.926.06'
$ws.Range("E15").Value = '  -5.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008527'
$ws.Range("E16").Value = '  -6.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.55'
$ws.Range("E17").Value = '  -12.94%  '
$ws.Range("D18").Value = '26.631.79'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.990'
$ws.Range("E19").Value = '  -8.77%  '
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.98'
$ws.Range("E21").Value = '  -4.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '186.68'
$ws.Range("E22").Value = '  -11.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.263'
$ws.Range("E23").Value = '  -8.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.006'
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.66'
$ws.Range("E25").Value = '  -5.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.474'
$ws.Range("E26").Value = '  -7.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1164'
$ws.Range("E27").Value = '  -7.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.88'
$ws.Range("E28").Value = '  -3.52%  '
$ws.Range("E29").Value = '  -6.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05752'
$ws.Range("E30").Value = '  -6.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.343'
$ws.Range("E31").Value = '  -5.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.521'
$ws.Range("E32").Value = '  -7.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.513'
$ws.Range("E33").Value = '  -7.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.647'
$ws.Range("E34").Value = '  -5.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.021'
$ws.Range("E35").Value = '  -2.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6000'
$ws.Range("E36").Value = '  -6.47%  '
$ws.Range("E37").Value = '  -5.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.685'
$ws.Range("E38").Value = '  -1.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01623'
$ws.Range("E39").Value = '  -4.46%  '
$ws.Range("D40").Value = '1.103.55'
$ws.Range("E40").Value = '  -3.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8578'
$ws.Range("E41").Value = '  -2.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.838'
$ws.Range("E42").Value = '  -9.01%  '
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.31'
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").Value = '1.850.76'
$ws.Range("E45").Value = '  -5.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000117'
$ws.Range("E46").Value = '  +5.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.67'
$ws.Range("E47").Value = '  -5.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("E48").Value = '  +0.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.110'
$ws.Range("E49").Value = '  -2.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4324'
$ws.Range("E50").Value = '  -3.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05232'
$ws.Range("E51").Value = '  -4.71%  '
